# Update Mappings 22 Ontologies
# Adds a new "REX_DEF" column (F) to the mapping sheet, populated with the
# default value "[]" for every existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in F1, matching the formatting already used by the other
# header cells (B1:E1).
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("F1").Value = "REX_DEF"

# Populate the new column for every data row with the placeholder "[]".
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 6).Value = "[]"
}
